$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligand average / total expression value (same for all data rows 2-5)
$ws.Range("G2:G5").Value = 0.06624833333333334
$ws.Range("H2:H5").Value = 0.198745

# Row 2 (MuSCs -> ECs)
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 0.12656200847
$ws.Range("R2").Value = 1.13905807623
$ws.Range("S2").Value = 0.01809124304049503
$ws.Range("T2").Value = 0.01809124304049503

# Row 3 (MuSCs -> FAPs) - M3/N3 unchanged
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 2.119632782063333
$ws.Range("R3").Value = 19.07669503857
$ws.Range("S3").Value = 0.302988173785169
$ws.Range("T3").Value = 0.302988173785169

# Row 4 (MuSCs -> MuSCs)
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 2.508075843415
$ws.Range("R4").Value = 22.572682590735
$ws.Range("S4").Value = 0.3585136661130873
$ws.Range("T4").Value = 0.3585136661130873

# Row 5 (MuSCs -> Resolving-Mac)
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 2.241490143058889
$ws.Range("R5").Value = 20.17341128753
$ws.Range("S5").Value = 0.3204069170612486
$ws.Range("T5").Value = 0.3204069170612486
